# ---------------------------------------------------------------------------
# UniformA-HW10.xlsx: "New simulation files for schemes report"
#
# The workbook lists averaging/HKL schemes (rows, col B) against a set of
# Miller-index / pairing columns (row 2, cols C:W). This edit:
#   1. Drops the now-unused HKL columns X:AG (the sheet only needs C:W).
#   2. Re-labels row 2 (C2:W2) because the HKL/pairing label order changed.
#   3. Renames the existing scheme rows (3-19) to the new scheme names.
#   4. Appends 10 new scheme rows (20-29) for the newly added simulations,
#      each flagged "1" across every HKL/pairing column, same as the rest.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the trailing HKL columns X:AG (rows 1-2 only used them)
$ws.Range("X1:AG2").Clear()

# 2) Re-point the HKL header row (row 2, C2:W2) at the reordered set of
#    Miller-index / pairing labels
$hklCols = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W")
$hklLabels = @(
    "[4, 2, 0]",
    "[4, 0, 0]",
    "[2, 0, 0]",
    "[2, 2, 0]",
    "[3, 3, 3]",
    "[1, 1, 1]",
    "[2, 2, 2]",
    "[3, 3, 1]",
    "[3, 1, 1]",
    "[4, 2, 2]",
    "[5, 1, 1]",
    "1Pair-A",
    "1Pair-B",
    "2Pairs-A",
    "2Pairs-B",
    "3Pairs-A",
    "3Pairs-B",
    "3Pairs-C",
    "4Pairs",
    "5A4F",
    "MaxUnique"
)
for ($i = 0; $i -lt $hklCols.Length; $i++) {
    $ws.Range("$($hklCols[$i])2").Value = $hklLabels[$i]
}

# 3) Relabel the existing scheme rows (3-19): the scheme-name list was
#    reordered/renamed; the per-column flags (C:W) are unchanged.
$existingSchemeNames = @(
    "Spiral5",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Thomas Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD"
)
for ($i = 0; $i -lt $existingSchemeNames.Length; $i++) {
    $row = 3 + $i
    $ws.Range("B$row").Value = $existingSchemeNames[$i]
}

# 4) Append the new simulation-scheme rows (20-29), mirroring the layout of
#    rows 3-19: col A = scheme index, col B = scheme name, C:W = 1 (flags).
#    Copy row-by-row - a single copy into a taller destination only fills
#    the first row on this host, so the source is re-copied per new row.
$newSchemeNames = @(
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex"
)
for ($i = 0; $i -lt $newSchemeNames.Length; $i++) {
    $row = 20 + $i
    $ws.Range("A19:W19").Copy($ws.Range("A$($row):W$($row)"))
    $ws.Range("A$row").Value = 18 + $i
    $ws.Range("B$row").Value = $newSchemeNames[$i]
}
